{"js": "// 1) Title block: \"CSE3215 DIGITAL LOGIC DESIGN PROJECT\" -> \"CSE3215 DIGITAL LOGIC DESIGN TERM PROJECT\"\n//    (insert \"TERM \" right before \"PROJECT\", keeping the existing bold/sz32 formatting)\n//    and remove the following \"Iteration 1: Assembly Language\" paragraph entirely.\nconst body = context.document.body;\n\nconst titleResults = body.search(\"CSE3215 DIGITAL LOGIC DESIGN PROJECT\", { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  const titleRange = titleResults.items[0];\n  const projectResults = titleRange.search(\"PROJECT\", { matchCase: true });\n  projectResults.load(\"items\");\n  await context.sync();\n\n  if (projectResults.items.length > 0) {\n    projectResults.items[0].insertText(\"TERM \", Word.InsertLocation.before);\n    await context.sync();\n  }\n}\n\nconst iterationResults = body.search(\"Iteration 1: Assembly Language\", { matchCase: true });\niterationResults.load(\"items\");\nawait context.sync();\n\nif (iterationResults.items.length > 0) {\n  const iterRange = iterationResults.items[0];\n  const iterParagraph = iterRange.paragraphs.getFirst();\n  iterParagraph.delete();\n  await context.sync();\n}\n\n// 2) CMP row of the instruction-set table: swap the \"OP1\" / \"0000\" cell values\n//    (3rd and 4th data cells of the CMP row).\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (let t = 0; t < tables.items.length; t++) {\n  const table = tables.items[t];\n  table.load(\"values\");\n  await context.sync();\n\n  for (let r = 0; r < table.values.length; r++) {\n    const rowValues = table.values[r];\n    if (rowValues[0] === \"CMP\" && rowValues[2] === \"OP1\" && rowValues[3] === \"0000\") {\n      const cellOp1 = table.getCell(r, 2);\n      const cellZero = table.getCell(r, 3);\n      const rangeOp1 = cellOp1.body.paragraphs.getFirst().getRange();\n      const rangeZero = cellZero.body.paragraphs.getFirst().getRange();\n      rangeOp1.insertText(\"0000\", Word.InsertLocation.replace);\n      rangeZero.insertText(\"OP1\", Word.InsertLocation.replace);\n      await context.sync();\n    }\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Title block: \"CSE3215 DIGITAL LOGIC DESIGN PROJECT\" -> \"CSE3215 DIGITAL LOGIC DESIGN TERM PROJECT\"\n#    Insert \"TERM \" right before \"PROJECT\" (keeps the existing bold/sz32 formatting),\n#    then remove the following \"Iteration 1: Assembly Language\" paragraph entirely.\n$titleRng = $d.Content\n$titleRng.Find.Execute(\"PROJECT\") | Out-Null\n$titleRng.Collapse(1) | Out-Null\n$titleRng.InsertBefore(\"TERM \")\n\n$iterRng = $d.Content\n$iterRng.Find.Execute(\"Iteration 1: Assembly Language\") | Out-Null\n$iterRng.Expand(4) | Out-Null\n$iterRng.Delete()\n\n# 2) CMP row of the instruction-set table: swap the \"OP1\" / \"0000\" cell values\n#    (3rd and 4th data cells of the CMP row).\nforeach ($tbl in $d.Tables) {\n  for ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    $firstCellText = $tbl.Cell($r, 1).Range.Text -replace [char]13, '' -replace [char]7, ''\n    if ($firstCellText -eq \"CMP\") {\n      $cellOp1 = $tbl.Cell($r, 3)\n      $cellZero = $tbl.Cell($r, 4)\n      $op1Text = $cellOp1.Range.Text -replace [char]13, '' -replace [char]7, ''\n      $zeroText = $cellZero.Range.Text -replace [char]13, '' -replace [char]7, ''\n      if ($op1Text -eq \"OP1\" -and $zeroText -eq \"0000\") {\n        $cellOp1.Range.Text = \"0000\"\n        $cellZero.Range.Text = \"OP1\"\n      }\n    }\n  }\n}\n"}
